# SteadyCalibrationCurve.xlsx — switch the calibration curve's independent
# variable (column A, "Known"/reference values on Sheet2) from an 80-unit
# step series (0, 80, 160, ... 3840) to a 75-unit step series
# (0, 75, 150, ... 3600), reflecting newly captured automatic-detection /
# sample-measurement data. Columns B (raw sensor reading) and C (state) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet2's used range is A1:C49 -- rewrite the 49 values of column A in
# place, row by row, so dependent chart series (categories = Sheet2!$A:$A)
# track the new spacing.
for ($row = 1; $row -le 49; $row++) {
    $ws.Cells.Item($row, 1).Value = ($row - 1) * 75
}

# Restore the cursor/selection left behind by the edit (the workbook was
# last saved with E38 selected, scrolled so row 23 is at the top).
[void]$ws.Range("E38").Select()
